$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 (E17/F17): change "TO DO" -> "OK".
# Reuse the formatting already used for "OK" cells elsewhere on the sheet
# (e.g. E2:F2) via copy/paste-special so the engine reuses the existing
# cellXfs style entry instead of minting a brand-new one.
$ws.Range("E2:F2").Copy()
$ws.Range("E17:F17").PasteSpecial(-4122)
$ws.Range("E17:F17").Value = "OK"

# Update the active selection to C17.
$ws.Range("C17").Select()
